$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 38: SQL entry about parsing string/numeric values into a date.
# Shared strings must be created in the same order as the target file:
# the long SQL snippet (index 69) before the short label (index 70).
$ws.Range("C38").Value = " select YEAR(TIMESTAMP_FORMAT(DIGITS(BDPPOSTD),'YYYYMMDD')) `n from iprod/bdppay                                          "
$ws.Range("B38").Value = "Parse string/numeric to date"
$ws.Range("A38").Value = "SQL"

# Match the formatting of the other multi-line SQL rows (wrapped text cells).
$ws.Range("B38").WrapText = $true
$ws.Range("C38").WrapText = $true

# Row height matches the rest of the sheet.
$ws.Rows.Item(38).RowHeight = 37.5

# Update the selection to the new entry point, like the author did after typing the row.
$ws.Range("B39").Select()
